$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; this shifts the former rows 41-123 down to 42-124
# and also extends the sheet dimension (A1:T123 -> A1:T124) automatically.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new observation.
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 45272
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100103
$ws.Range("H41").Value = "Frutos de hueso (carozo)"
$ws.Range("I41").Value = 100103006
$ws.Range("J41").Value = "Nectarín"
$ws.Range("K41").Value = "Early Diamond"
$ws.Range("L41").Value = "Segunda"
$ws.Range("M41").Value = 300
$ws.Range("N41").Value = 20000
$ws.Range("O41").Value = 25000
$ws.Range("P41").Value = 22500
$ws.Range("Q41").Value = "$/bandeja 18 kilos granel"
$ws.Range("R41").Value = "Región de O'Higgins"
$ws.Range("S41").Value = 1250
$ws.Range("T41").Value = 18
